# Auto-generated cell value updates derived from the OOXML diff.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H42").Value = 113.8
$ws1.Range("J42").Value = 115.333336
$ws1.Range("L42").Value = 346.000008
$ws1.Range("N42").Value = -806.000008
$ws1.Range("H57").Value = 22042
$ws1.Range("J57").Value = 22042
$ws1.Range("L57").Value = 66126
$ws1.Range("N57").Value = -67124
$ws1.Range("H113").Value = 45458216
$ws1.Range("I113").Value = 76926150
$ws1.Range("J113").Value = 4532
$ws1.Range("K113").Value = 76926150
$ws1.Range("L113").Value = 4532
$ws1.Range("M113").Value = -76922896
$ws1.Range("N113").Value = -11040
$ws1.Range("H129").Value = 294893.7
$ws1.Range("J129").Value = 334162.84
$ws1.Range("L129").Value = 1002488.52
$ws1.Range("N129").Value = -1012488.52
$ws1.Range("H137").Value = 52763.65
$ws1.Range("I137").Value = 3206
$ws1.Range("J137").Value = 127100.125
$ws1.Range("K137").Value = 9618
$ws1.Range("L137").Value = 381300.375
$ws1.Range("M137").Value = -7068
$ws1.Range("N137").Value = -386400.375

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H2").Value = 1343.5938
$ws2.Range("I2").Value = 1163.875
$ws2.Range("J2").Value = 1882.75
$ws2.Range("K2").Value = 1163.875
$ws2.Range("L2").Value = 1882.75
$ws2.Range("M2").Value = -1050.875
$ws2.Range("N2").Value = -2108.75
$ws2.Range("H61").Value = 3097.9524
$ws2.Range("I61").Value = 2526.4614
$ws2.Range("J61").Value = 4026.625
$ws2.Range("K61").Value = 2526.4614
$ws2.Range("L61").Value = 4026.625
$ws2.Range("M61").Value = -2314.4614
$ws2.Range("N61").Value = -4450.625
$ws2.Range("H74").Value = 66667540
$ws2.Range("I74").Value = 142857570
$ws2.Range("J74").Value = 1258.625
$ws2.Range("K74").Value = 142857570
$ws2.Range("L74").Value = 1258.625
$ws2.Range("M74").Value = -142856696
$ws2.Range("N74").Value = -3006.625
$ws2.Range("H77").Value = 66667540
$ws2.Range("I77").Value = 142857570
$ws2.Range("J77").Value = 1258.625
$ws2.Range("K77").Value = 714287850
$ws2.Range("L77").Value = 6293.125
$ws2.Range("M77").Value = -714283482
$ws2.Range("N77").Value = -15029.125
$ws2.Range("H116").Value = 1343.5938
$ws2.Range("I116").Value = 1163.875
$ws2.Range("J116").Value = 1882.75
$ws2.Range("K116").Value = 1163.875
$ws2.Range("L116").Value = 1882.75
$ws2.Range("M116").Value = 1130.125
$ws2.Range("N116").Value = -6470.75
$ws2.Range("H132").Value = 14024.738
$ws2.Range("I132").Value = 1841
$ws2.Range("J132").Value = 41203.848
$ws2.Range("K132").Value = 5523
$ws2.Range("L132").Value = 123611.544
$ws2.Range("M132").Value = -2993
$ws2.Range("N132").Value = -128671.544
$ws2.Range("H134").Value = 47282.6
$ws2.Range("J134").Value = 47282.6
$ws2.Range("L134").Value = 47282.6
$ws2.Range("N134").Value = -57422.6
$ws2.Range("H135").Value = 37202
$ws2.Range("J135").Value = 37202
$ws2.Range("L135").Value = 37202
$ws2.Range("N135").Value = -47342
$ws2.Range("H136").Value = 3097.9524
$ws2.Range("I136").Value = 2526.4614
$ws2.Range("J136").Value = 4026.625
$ws2.Range("K136").Value = 7579.3842
$ws2.Range("L136").Value = 12079.875
$ws2.Range("M136").Value = -5029.3842
$ws2.Range("N136").Value = -17179.875

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H3").Value = 1343.5938
$ws3.Range("I3").Value = 1163.875
$ws3.Range("J3").Value = 1882.75
$ws3.Range("K3").Value = 1163.875
$ws3.Range("L3").Value = 1882.75
$ws3.Range("M3").Value = -1049.875
$ws3.Range("N3").Value = -2110.75
$ws3.Range("H105").Value = 3426.5264
$ws3.Range("I105").Value = 3317.8333
$ws3.Range("K105").Value = 3317.8333
$ws3.Range("M105").Value = -1570.8333
$ws3.Range("H134").Value = 33268.91
$ws3.Range("I134").Value = 42889.08
$ws3.Range("J134").Value = 2003.375
$ws3.Range("K134").Value = 128667.24
$ws3.Range("L134").Value = 6010.125
$ws3.Range("M134").Value = -126132.24
$ws3.Range("N134").Value = -11080.125

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H31").Value = 9991.512000000001
$ws4.Range("I31").Value = 14561.417
$ws4.Range("K31").Value = 14561.417
$ws4.Range("M31").Value = -14266.417
$ws4.Range("H34").Value = 9991.512000000001
$ws4.Range("I34").Value = 14561.417
$ws4.Range("K34").Value = 14561.417
$ws4.Range("M34").Value = -14359.417
$ws4.Range("H58").Value = 11542.617
$ws4.Range("I58").Value = 918.9143
$ws4.Range("K58").Value = 918.9143
$ws4.Range("M58").Value = -715.9143
$ws4.Range("H132").Value = 33953.824
$ws4.Range("I132").Value = 45009
$ws4.Range("J132").Value = 7421.4
$ws4.Range("K132").Value = 135027
$ws4.Range("L132").Value = 22264.2
$ws4.Range("M132").Value = -132497
$ws4.Range("N132").Value = -27324.2
$ws4.Range("H134").Value = 1284.4286
$ws4.Range("I134").Value = 963.4286
$ws4.Range("K134").Value = 2890.2858
$ws4.Range("M134").Value = -355.2857999999997
$ws4.Range("H136").Value = 11542.617
$ws4.Range("I136").Value = 918.9143
$ws4.Range("K136").Value = 2756.7429
$ws4.Range("M136").Value = -206.7429000000002

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H43").Value = 0
$ws5.Range("J43").Value = 0
$ws5.Range("L43").Value = 0
$ws5.Range("N43").Value = ""
$ws5.Range("H131").Value = 766.3200000000001
$ws5.Range("J131").Value = 770.4433
$ws5.Range("L131").Value = 2311.3299
$ws5.Range("N131").Value = -12391.3299

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H80").Value = 3526.1
$ws6.Range("I80").Value = 3115.4285
$ws6.Range("K80").Value = 3115.4285
$ws6.Range("M80").Value = -2117.4285
$ws6.Range("H83").Value = 3526.1
$ws6.Range("I83").Value = 3115.4285
$ws6.Range("K83").Value = 15577.1425
$ws6.Range("M83").Value = -10585.1425
$ws6.Range("H132").Value = 49159.637
$ws6.Range("I132").Value = 54455.8
$ws6.Range("J132").Value = 41011.69
$ws6.Range("K132").Value = 163367.4
$ws6.Range("L132").Value = 123035.07
$ws6.Range("M132").Value = -160837.4
$ws6.Range("N132").Value = -128095.07

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H7").Value = 5729.647
$ws7.Range("J7").Value = 3337.5
$ws7.Range("L7").Value = 3337.5
$ws7.Range("N7").Value = -3561.5
$ws7.Range("H16").Value = 470.85715
$ws7.Range("I16").Value = 469.2
$ws7.Range("J16").Value = 475
$ws7.Range("K16").Value = 469.2
$ws7.Range("L16").Value = 475
$ws7.Range("M16").Value = -299.2
$ws7.Range("N16").Value = -815
$ws7.Range("H22").Value = 1672
$ws7.Range("J22").Value = 1958.75
$ws7.Range("L22").Value = 1958.75
$ws7.Range("N22").Value = -2548.75
$ws7.Range("H27").Value = 1672
$ws7.Range("J27").Value = 1958.75
$ws7.Range("L27").Value = 1958.75
$ws7.Range("N27").Value = -2172.75
$ws7.Range("H126").Value = 5729.647
$ws7.Range("J126").Value = 3337.5
$ws7.Range("L126").Value = 10012.5
$ws7.Range("N126").Value = -14952.5
$ws7.Range("H132").Value = 3005.7334
$ws7.Range("I132").Value = 1908.8
$ws7.Range("J132").Value = 5199.6
$ws7.Range("K132").Value = 5726.4
$ws7.Range("L132").Value = 15598.8
$ws7.Range("M132").Value = -3196.4
$ws7.Range("N132").Value = -20658.8
$ws7.Range("H133").Value = 27400
$ws7.Range("J133").Value = 27400
$ws7.Range("L133").Value = 27400
$ws7.Range("N133").Value = -32460
$ws7.Range("H136").Value = 28158
$ws7.Range("I136").Value = 42833.5
$ws7.Range("K136").Value = 128500.5
$ws7.Range("M136").Value = -125950.5

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H132").Value = 2107.647
$ws8.Range("I132").Value = 1411.4445
$ws8.Range("J132").Value = 2890.875
$ws8.Range("K132").Value = 4234.333500000001
$ws8.Range("L132").Value = 8672.625
$ws8.Range("M132").Value = -1704.333500000001
$ws8.Range("N132").Value = -13732.625
$ws8.Range("H136").Value = 33335464
$ws8.Range("I136").Value = 62502110
$ws8.Range("J136").Value = 2150.6428
$ws8.Range("K136").Value = 187506330
$ws8.Range("L136").Value = 6451.928400000001
$ws8.Range("M136").Value = -187503780
$ws8.Range("N136").Value = -11551.9284

